$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03835166666666667
$ws.Range("H2").Value = 0.115055
$ws.Range("I2").Value = 0.0181239951898282
$ws.Range("J2").Value = 0.0181239951898282
$ws.Range("M2").Value = 0.01031333333333333
$ws.Range("N2").Value = 0.03094
$ws.Range("O2").Value = 0.001146416507271297
$ws.Range("P2").Value = 0.001146416507271297
$ws.Range("Q2").Value = 0.0003955335222222222
$ws.Range("R2").Value = 0.0035598017
$ws.Range("S2").Value = 0.00002077764726332462
$ws.Range("T2").Value = 0.00002077764726332462

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03835166666666667
$ws.Range("H3").Value = 0.115055
$ws.Range("I3").Value = 0.0181239951898282
$ws.Range("J3").Value = 0.0181239951898282
$ws.Range("O3").Value = 0.8303652607489888
$ws.Range("P3").Value = 0.8303652607489886
$ws.Range("Q3").Value = 0.2864903760822222
$ws.Range("R3").Value = 2.57841338474
$ws.Range("S3").Value = 0.01504953599161511
$ws.Range("T3").Value = 0.01504953599161511

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03835166666666667
$ws.Range("H4").Value = 0.115055
$ws.Range("I4").Value = 0.0181239951898282
$ws.Range("J4").Value = 0.0181239951898282
$ws.Range("M4").Value = 1.515746
$ws.Range("N4").Value = 4.547238
$ws.Range("O4").Value = 0.16848832274374
$ws.Range("P4").Value = 0.16848832274374
$ws.Range("Q4").Value = 0.05813138534333333
$ws.Range("R4").Value = 0.52318246809
$ws.Range("S4").Value = 0.003053681550949765
$ws.Range("T4").Value = 0.003053681550949765

$ws.Range("I5").Value = 0.3727881574250648
$ws.Range("J5").Value = 0.3727881574250648
$ws.Range("M5").Value = 0.01031333333333333
$ws.Range("N5").Value = 0.03094
$ws.Range("O5").Value = 0.001146416507271297
$ws.Range("P5").Value = 0.001146416507271297
$ws.Range("Q5").Value = 0.008135635184444443
$ws.Range("R5").Value = 0.07322071665999999
$ws.Range("S5").Value = 0.0004273704973873451
$ws.Range("T5").Value = 0.000427370497387345

$ws.Range("I6").Value = 0.3727881574250648
$ws.Range("J6").Value = 0.3727881574250648
$ws.Range("O6").Value = 0.8303652607489888
$ws.Range("P6").Value = 0.8303652607489886
$ws.Range("S6").Value = 0.309550335544399
$ws.Range("T6").Value = 0.3095503355443989

$ws.Range("I7").Value = 0.3727881574250648
$ws.Range("J7").Value = 0.3727881574250648
$ws.Range("M7").Value = 1.515746
$ws.Range("N7").Value = 4.547238
$ws.Range("O7").Value = 0.16848832274374
$ws.Range("P7").Value = 0.16848832274374
$ws.Range("Q7").Value = 1.195690674364667
$ws.Range("R7").Value = 10.761216069282
$ws.Range("S7").Value = 0.06281045138327848
$ws.Range("T7").Value = 0.06281045138327848

$ws.Range("G8").Value = 1.288873333333333
$ws.Range("H8").Value = 3.86662
$ws.Range("I8").Value = 0.6090878473851071
$ws.Range("J8").Value = 0.609087847385107
$ws.Range("M8").Value = 0.01031333333333333
$ws.Range("N8").Value = 0.03094
$ws.Range("O8").Value = 0.001146416507271297
$ws.Range("P8").Value = 0.001146416507271297
$ws.Range("Q8").Value = 0.01329258031111111
$ws.Range("R8").Value = 0.1196332228
$ws.Range("S8").Value = 0.0006982683626206271
$ws.Range("T8").Value = 0.0006982683626206269

$ws.Range("G9").Value = 1.288873333333333
$ws.Range("H9").Value = 3.86662
$ws.Range("I9").Value = 0.6090878473851071
$ws.Range("J9").Value = 0.609087847385107
$ws.Range("O9").Value = 0.8303652607489888
$ws.Range("P9").Value = 0.8303652607489886
$ws.Range("Q9").Value = 9.627998939351112
$ws.Range("R9").Value = 86.65199045416
$ws.Range("S9").Value = 0.5057653892129748
$ws.Range("T9").Value = 0.5057653892129745

$ws.Range("G10").Value = 1.288873333333333
$ws.Range("H10").Value = 3.86662
$ws.Range("I10").Value = 0.6090878473851071
$ws.Range("J10").Value = 0.609087847385107
$ws.Range("M10").Value = 1.515746
$ws.Range("N10").Value = 4.547238
$ws.Range("O10").Value = 0.16848832274374
$ws.Range("P10").Value = 0.16848832274374
$ws.Range("Q10").Value = 1.953604599506667
$ws.Range("R10").Value = 17.58244139556
$ws.Range("S10").Value = 0.1026241898095118
$ws.Range("T10").Value = 0.1026241898095118

